$wb = $excel.ActiveWorkbook

# --- Sheet "Our Results": update BM25_MonoT5 row (row 3) ---
$wsResults = $wb.Worksheets.Item("Our Results")
$wsResults.Range("C3").Value = 0.4925
$wsResults.Range("D3").Value = 0.518
$wsResults.Range("E3").Value = 0.5303
$wsResults.Range("F3").Value = 0.4836
$wsResults.Range("G3").Value = 0.6742
$wsResults.Range("H3").Value = 0.7497

# --- Sheet "BEIR Comparison": update BM25_MonoT5 row (row 3) ---
$wsComparison = $wb.Worksheets.Item("BEIR Comparison")
$wsComparison.Range("C3").Value = 0.518

# F3 holds a text label like "-2.8%" (not a numeric percentage), so force
# text entry (leading apostrophe) to keep it a literal string instead of
# having it auto-converted into a percentage number, then restore the
# bold black font used for "≈ Close" deltas.
$cellF3 = $wsComparison.Range("F3")
$cellF3.Value = "'-2.8%"
$cellF3.Font.Bold = $true
